$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1979949874686717
$ws.Range("C2").Value = 0.5338345864661654
$ws.Range("J2").Value = 0.012531328320802
$ws.Range("P2").Value = 0.1528822055137845
$ws.Range("S2").Value = 0.1027568922305764
$ws.Range("B3").Value = 0.01327433628318584
$ws.Range("C3").Value = 0.06194690265486726
$ws.Range("J3").Value = 0.03539823008849557
$ws.Range("P3").Value = 0.7256637168141593
$ws.Range("S3").Value = 0.163716814159292
$ws.Range("J4").Value = 0.07547169811320754
$ws.Range("P4").Value = 0.6037735849056604
$ws.Range("S4").Value = 0.3207547169811321
$ws.Range("B6").Value = 0.06439393939393939
$ws.Range("D6").Value = 0.01515151515151515
$ws.Range("F6").Value = 0.04924242424242424
$ws.Range("J6").Value = 0.2689393939393939
$ws.Range("O6").Value = 0.01136363636363636
$ws.Range("Q6").Value = 0.2121212121212121
$ws.Range("S6").Value = 0.3484848484848485
$ws.Range("B7").Value = 0.1627906976744186
$ws.Range("D7").Value = 0.009302325581395349
$ws.Range("F7").Value = 0.05116279069767442
$ws.Range("J7").Value = 0.1255813953488372
$ws.Range("O7").Value = 0.02325581395348837
$ws.Range("Q7").Value = 0.2093023255813954
$ws.Range("R7").Value = 0.06511627906976744
$ws.Range("S7").Value = 0.3534883720930233
$ws.Range("B8").Value = 0.1126482213438735
$ws.Range("D8").Value = 0.02371541501976284
$ws.Range("E8").Value = 0.001976284584980237
$ws.Range("F8").Value = 0.07509881422924901
$ws.Range("J8").Value = 0.09683794466403162
$ws.Range("O8").Value = 0.02371541501976284
$ws.Range("Q8").Value = 0.1660079051383399
$ws.Range("R8").Value = 0.06719367588932806
$ws.Range("S8").Value = 0.4328063241106719
$ws.Range("B9").Value = 0.1194690265486726
$ws.Range("D9").Value = 0.008849557522123894
$ws.Range("F9").Value = 0.07079646017699115
$ws.Range("J9").Value = 0.08849557522123894
$ws.Range("O9").Value = 0.02654867256637168
$ws.Range("Q9").Value = 0.2035398230088496
$ws.Range("R9").Value = 0.1106194690265487
$ws.Range("S9").Value = 0.3716814159292036
$ws.Range("B10").Value = 0.1182572614107884
$ws.Range("D10").Value = 0.02351313969571231
$ws.Range("E10").Value = 0.0006915629322268327
$ws.Range("F10").Value = 0.0698478561549101
$ws.Range("J10").Value = 0.1078838174273859
$ws.Range("O10").Value = 0.01936376210235131
$ws.Range("Q10").Value = 0.1922544951590595
$ws.Range("R10").Value = 0.08437067773167359
$ws.Range("S10").Value = 0.3838174273858921
$ws.Range("G11").Value = 0.15
$ws.Range("J11").Value = 0.1375
$ws.Range("K11").Value = 0.2325
$ws.Range("L11").Value = 0.47
$ws.Range("S11").Value = 0.01
$ws.Range("F12").Value = 0.005208333333333333
$ws.Range("G12").Value = 0.6302083333333334
$ws.Range("J12").Value = 0.2447916666666667
$ws.Range("K12").Value = 0.02083333333333333
$ws.Range("L12").Value = 0.01041666666666667
$ws.Range("S12").Value = 0.08854166666666667
$ws.Range("G13").Value = 0.6101694915254238
$ws.Range("J13").Value = 0.3220338983050847
$ws.Range("S13").Value = 0.06779661016949153
$ws.Range("F15").Value = 0.02232142857142857
$ws.Range("H15").Value = 0.1339285714285714
$ws.Range("I15").Value = 0.08482142857142858
$ws.Range("J15").Value = 0.3348214285714285
$ws.Range("K15").Value = 0.05357142857142857
$ws.Range("M15").Value = 0.008928571428571428
$ws.Range("O15").Value = 0.0625
$ws.Range("S15").Value = 0.2991071428571428
$ws.Range("F16").Value = 0.02
$ws.Range("H16").Value = 0.168
$ws.Range("I16").Value = 0.07199999999999999
$ws.Range("J16").Value = 0.42
$ws.Range("K16").Value = 0.112
$ws.Range("M16").Value = 0.016
$ws.Range("O16").Value = 0.032
$ws.Range("S16").Value = 0.16
$ws.Range("F17").Value = 0.02514506769825919
$ws.Range("H17").Value = 0.1760154738878143
$ws.Range("I17").Value = 0.0735009671179884
$ws.Range("J17").Value = 0.3945841392649903
$ws.Range("K17").Value = 0.1160541586073501
$ws.Range("M17").Value = 0.01547388781431335
$ws.Range("O17").Value = 0.06189555125725339
$ws.Range("S17").Value = 0.137330754352031
$ws.Range("F18").Value = 0.0154639175257732
$ws.Range("H18").Value = 0.1649484536082474
$ws.Range("I18").Value = 0.09793814432989691
$ws.Range("J18").Value = 0.4072164948453608
$ws.Range("K18").Value = 0.09793814432989691
$ws.Range("M18").Value = 0.02061855670103093
$ws.Range("O18").Value = 0.06701030927835051
$ws.Range("S18").Value = 0.1288659793814433
$ws.Range("F19").Value = 0.01806020066889632
$ws.Range("H19").Value = 0.2147157190635451
$ws.Range("I19").Value = 0.0882943143812709
$ws.Range("J19").Value = 0.3625418060200669
$ws.Range("K19").Value = 0.1204013377926421
$ws.Range("M19").Value = 0.02809364548494983
$ws.Range("N19").Value = 0.002675585284280936
$ws.Range("O19").Value = 0.05217391304347826
$ws.Range("S19").Value = 0.1130434782608696
